$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 957
$wsExhibit.Range("F3").Value = 1863
$wsExhibit.Range("F4").Value = 416

# Sheet "全部类型": same events appear again, shifted down two rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 957
$wsAll.Range("F5").Value = 1863
$wsAll.Range("F6").Value = 416
